$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8:83 down to 9:84
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with its data
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44817
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = 300000001
$ws.Cells.Item(8, 7).Value = "Rabanito"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 55
$ws.Cells.Item(8, 11).Value = 8000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 8000
$ws.Cells.Item(8, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(8, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(8, 16).Value = 667
$ws.Cells.Item(8, 17).Value = 12
$ws.Cells.Item(8, 18).Value = "Hortaliza"
